$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2: change FLOOR formulas to CEILING in columns D, E, I (rows 2-37) ---
$ws2.Range("D2:D37").Formula = "=CEILING(B2/25, 1)"
$ws2.Range("E2:E37").Formula = "=CEILING(B2/30, 1)"
$ws2.Range("I2:I37").Formula = "=CEILING(C2/25, 1)"

# --- Sheet2: change column M values (rows 2-37) to 11 ---
$ws2.Range("M2:M37").Value = 11

# --- Sheet2: view state - topLeftCell D1, selection N20 ---
$ws2.Activate()
$ws2.Application.ActiveWindow.ScrollColumn = 4
$ws2.Range("N20").Select()

# --- Sheet1: selection K5 ---
$ws1.Activate()
$ws1.Range("K5").Select()

# Re-activate Sheet2 as the final active tab (tabSelected="1" stays on Sheet2)
$ws2.Activate()
